$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old)
    if (-not $ok) {
        throw "Find failed for: $old"
    }
    $rng.Text = $new
}

# 1. Title
Replace-Text "Kepler's Laws: Orchestrating the Cosmic Dance" "The Celestial Symphony: Unraveling the Enigma of the Cosmos"

# 2. Author
Replace-Text "Amelia Stevens" "Alicia White"

# 3. Email: "newtonphysics@researchhub" (+ existing ".org") -> "alicia" + "." + "white@validdomain" (+ existing ".org")
Replace-Text "newtonphysics@researchhub" "alicia.white@validdomain"

# 4. Body paragraph 1 sentences
Replace-Text "Through the vast expanse of the cosmos, planets and celestial bodies pirouette in a delicate balance, their movements governed by the timeless laws of Johannes Kepler" "In the vast expanse of the cosmos, there lies a symphony of celestial wonders, an intricate tapestry of stars, planets, galaxies, and cosmic phenomena"
Replace-Text " His pioneering insights into planetary motions revolutionized our understanding of the universe, unveiling the intricate harmony underpinning the celestial ballet" " Since time immemorial, humans have gazed upon the night sky with a mix of awe and curiosity, seeking to understand the enigmatic workings of the universe"
Replace-Text " Kepler's laws, borne from meticulous observations and mathematical brilliance, continue to enchant and inform our comprehension of the cosmos, providing a framework for comprehending the choreography of celestial bodies" " With the advent of modern science, we have embarked on a thrilling voyage of exploration, unraveling the secrets hidden within the cosmos"
Replace-Text "In the tapestry of astronomy, Kepler's name is etched in gold" "As we delve deeper into the celestial realm, we witness a breathtaking display of cosmic diversity"
Replace-Text " His meticulous observations of the heavens, driven by an insatiable curiosity and undeterred by the limitations of 17th-century technology, yielded groundbreaking discoveries" " From the fiery heart of our Sun to the distant reaches of far-off galaxies, each celestial body holds its unique story, a testament to the vastness and complexity of the universe"

# Merge: " He charted...elliptical paths" + "." + " Through painstaking...cosmic dance" -> one new sentence
Replace-Text " He charted the courses of planets, meticulously documenting their positions and velocities, unraveling the secrets of their elliptical paths. Through painstaking calculations and unwavering dedication, Kepler unveiled the mathematical harmonies governing planetary motion, orchestrating the cosmic dance" " The intricate dance of planets around their stars, the graceful ballet of moons orbiting their worlds, and the mesmerizing spectacle of stellar explosions paint a vibrant canvas of cosmic interaction"

Replace-Text "His laws, like sonorous melodies resonating through the cosmos, elucidated the intricate mechanisms that dictate the symphony of celestial bodies" "Unraveling the enigma of the cosmos is a testament to human ingenuity and our relentless pursuit of knowledge"
Replace-Text " His first law, like a conductor's precise baton, defines the elliptical paths of planets, guiding their graceful dance around the central sun" " Through observation, experimentation, and mathematical modeling, scientists have pieced together the intricate puzzle of the universe, revealing its fundamental laws and illuminating its deepest mysteries"

# Merge: " His second law...proximity to the sun" + "." + " The third law...mathematical precision" -> one new sentence
Replace-Text " His second law, revealing the variation of a planet's speed along its orbit, captures the rhythm of their journey, their velocity harmonizing with their proximity to the sun. The third law, a symphony of proportions, unveils the elegant relationship between a planet's orbital period and its mean distance from the sun, a celestial waltz governed by mathematical precision" " From the elegant simplicity of Kepler's laws to the profound insights of Einstein's relativity, our understanding of the cosmos has undergone a remarkable transformation"

# 5. Summary paragraph
Replace-Text "Kepler's Laws, derived from meticulous observations and mathematical rigor, unravel the intricate patterns of planetary motion, illuminating the harmony of the cosmos" "Our journey into the celestial symphony has revealed the intricate beauty and awe-inspiring grandeur of the universe"
Replace-Text " His first law dictates elliptical paths, the second reveals varying speeds along the orbit, and the third establishes a proportional relationship between orbital period and mean distance from the sun" " From ancient stargazers to modern scientists, humanity's quest to comprehend the cosmos has been a testament to our innate curiosity and unwavering pursuit of knowledge"

# Merge: " " + "These laws...generations to come" (with lastRenderedPageBreak) -> one new sentence, drop the page break
Replace-Text " These laws, like musical notes composing a celestial symphony, provide a profound understanding of the universal ballet, inspiring awe and wonder for generations to come" " As we continue to explore the vast expanse of space, we stand on the threshold of even greater discoveries, poised to unveil the mysteries that still shroud the cosmic realm"

# 6. Append an empty paragraph at the end of the document body
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.Text = "`r"
